# Update "想去人数" (F column) values across the four sheets to match the
# regenerated data snapshot (gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 8094
$ws.Range("F8").Value = 7019
$ws.Range("F9").Value = 7019
$ws.Range("F11").Value = 532
$ws.Range("F12").Value = 484
$ws.Range("F16").Value = 306
$ws.Range("F21").Value = 11493
$ws.Range("F22").Value = 117
$ws.Range("F23").Value = 2224
$ws.Range("F25").Value = 3071
$ws.Range("F28").Value = 2647
$ws.Range("F29").Value = 100
$ws.Range("F31").Value = 273
$ws.Range("F34").Value = 2349
$ws.Range("F36").Value = 1598
$ws.Range("F38").Value = 93
$ws.Range("F39").Value = 5767
$ws.Range("F40").Value = 1776
$ws.Range("F41").Value = 1243
$ws.Range("F42").Value = 825
$ws.Range("F43").Value = 158
$ws.Range("F47").Value = 1064
$ws.Range("F48").Value = 1508
$ws.Range("F49").Value = 96
$ws.Range("F50").Value = 1127

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 42
$ws.Range("F19").Value = 916

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 212
$ws.Range("F3").Value = 350

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 212
$ws.Range("F5").Value = 350
$ws.Range("F8").Value = 8094
$ws.Range("F12").Value = 7019
$ws.Range("F14").Value = 532
$ws.Range("F15").Value = 484
$ws.Range("F18").Value = 306
$ws.Range("F24").Value = 11493
$ws.Range("F25").Value = 117
$ws.Range("F26").Value = 2224
$ws.Range("F27").Value = 2224
$ws.Range("F28").Value = 3071
$ws.Range("F29").Value = 2647
$ws.Range("F31").Value = 273
$ws.Range("F34").Value = 2349
$ws.Range("F36").Value = 1598
$ws.Range("F38").Value = 93
$ws.Range("F39").Value = 5767
$ws.Range("F41").Value = 1776
$ws.Range("F43").Value = 1243
$ws.Range("F44").Value = 825
$ws.Range("F45").Value = 158
$ws.Range("F48").Value = 1064
$ws.Range("F49").Value = 1508
$ws.Range("F50").Value = 96
$ws.Range("F51").Value = 1127
